$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.06854566666666667
$ws.Range("M2").Value = 2.341355666666667
$ws.Range("N2").Value = 7.024067000000001
$ws.Range("O2").Value = 0.03973512964576821
$ws.Range("P2").Value = 0.0397351296457682
$ws.Range("Q2").Value = 0.1604897850754445
$ws.Range("R2").Value = 1.444408065679
$ws.Range("S2").Value = 0.03973512964576821
$ws.Range("T2").Value = 0.0397351296457682

# Row 3
$ws.Range("G3").Value = 0.06854566666666667
$ws.Range("O3").Value = 0.5779093692199981
$ws.Range("P3").Value = 0.5779093692199981
$ws.Range("S3").Value = 0.5779093692199981
$ws.Range("T3").Value = 0.5779093692199981

# Row 4
$ws.Range("G4").Value = 0.06854566666666667
$ws.Range("O4").Value = 0.3823555011342337
$ws.Range("P4").Value = 0.3823555011342337
$ws.Range("Q4").Value = 1.544329985745556
$ws.Range("S4").Value = 0.3823555011342337
$ws.Range("T4").Value = 0.3823555011342337
